$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 25 de Abril de 2020 a las 00:22"

$data = @(
    @("Estados Unidos",917347,30905,93283,772199,14932,1631,51865),
    @("España",219764,6740,92355,104885,7705,367,22524),
    @("Italia",192994,3021,60498,106527,2173,420,25969),
    @("Francia",159828,1645,43493,94090,4870,389,22245),
    @("Alemania",154999,1870,106800,42439,2908,185,5760),
    @("Reino Unido",143464,5386,0,123614,1559,768,19506),
    @("Turquia",104912,3122,21737,80575,1790,109,2600),
    @("Iran",88194,1168,66599,16021,3121,93,5574),
    @("China",82804,6,77257,915,57,0,4632),
    @("Rusia",68622,5849,5568,62439,2300,60,615),
    @("Brasil",52995,3503,26573,22752,8318,357,3670),
    @("Canada",44856,2746,15469,27093,557,147,2294),
    @("Belgica",44293,1496,10122,27492,970,189,6679),
    @("Paises Bajos",36535,806,0,31996,963,112,4289),
    @("Suiza",28677,181,21000,6088,386,40,1589),
    @("India",24447,1408,5496,18171,0,59,780),
    @("Portugal",22797,444,1228,20715,188,34,854),
    @("Ecuador",22719,11536,1366,20777,127,16,576),
    @("Peru",21648,734,7496,13518,505,62,634),
    @("Irlanda",18184,577,9233,7937,142,220,1014),
    @("Suecia",17567,812,550,14865,547,131,2152),
    @("Arabia Saudita",15102,1172,2049,12926,93,6,127),
    @("Austria",15071,69,11872,2669,156,8,530),
    @("Israel",15058,255,6003,8861,137,2,194),
    @("Japon",12368,0,1494,10546,259,0,328),
    @("Chile",12306,494,6327,5805,408,6,174),
    @("Singapur",12075,897,924,11139,26,0,12),
    @("Pakistan",11940,883,2755,8932,111,18,253),
    @("Mexico",11633,1089,2627,7937,378,99,1069),
    @("Polonia",10892,381,1944,8454,160,40,494),
    @("Corea del Sur",10708,6,8501,1967,55,0,240),
    @("Rumania",10417,321,2817,7033,241,22,567),
    @("Emiratos Arabes Unidos",9281,525,1760,7457,1,8,64),
    @("Bielorrusia",8773,751,1120,7590,92,3,63),
    @("Catar",8525,761,809,7706,72,0,10),
    @("Indonesia",8211,436,1002,6520,0,42,689),
    @("Dinamarca",8210,137,5526,2281,69,9,403),
    @("Ucrania",7647,477,601,6853,45,6,193),
    @("Serbia",7483,207,1094,6245,96,5,144),
    @("Noruega",7444,43,32,7213,53,5,199),
    @("Chequia",7236,49,2371,4651,76,4,214),
    @("Filipinas",7192,211,762,5953,1,15,477),
    @("Australia",6675,8,5136,1460,43,4,79),
    @("Republica Dominicana",5749,206,763,4719,135,2,267),
    @("Malasia",5691,88,3663,1932,41,1,96),
    @("Panama",5166,174,271,4749,86,2,146),
    @("Colombia",4881,320,1003,3653,98,10,225),
    @("Banglades",4689,503,112,4446,1,4,131),
    @("Finlandia",4395,111,2500,1718,61,5,177),
    @("Sudafrica",4220,267,1473,2668,36,4,79),
    @("Egipto",4092,201,1075,2723,0,7,294),
    @("Marruecos",3758,190,486,3114,1,3,158),
    @("Luxemburgo",3695,30,3007,603,25,2,85),
    @("Argentina",3435,0,976,2292,123,2,167),
    @("Argelia",3127,120,1408,1304,40,8,415),
    @("Moldavia",3110,184,755,2271,212,4,84),
    @("Tailandia",2854,15,2490,314,61,0,50),
    @("Kuwait",2614,215,613,1986,60,1,15),
    @("Barein",2518,301,1113,1397,1,0,8),
    @("Grecia",2490,27,577,1783,48,5,130),
    @("Kazajistan",2416,127,604,1787,29,5,25),
    @("Hungria",2383,99,401,1732,61,11,250),
    @("Croacia",2009,28,982,976,19,1,51),
    @("Uzbekistan",1804,46,621,1175,8,1,8),
    @("Oman",1790,74,325,1455,3,1,10),
    @("Islandia",1789,0,1542,237,5,0,10),
    @("Irak",1708,31,1204,418,0,3,86),
    @("Estonia",1605,13,206,1353,6,1,46),
    @("Armenia",1596,73,728,841,10,3,27),
    @("Azerbaiyan",1592,44,1013,558,14,1,21),
    @("Nueva Zelanda",1456,5,1095,344,1,1,17),
    @("Camerun",1430,96,668,719,20,0,43),
    @("Bosnia y Herzegovina",1421,8,538,828,4,1,55),
    @("Lituania",1410,12,430,940,17,0,40),
    @("Eslovenia",1373,7,211,1082,23,1,80),
    @("Eslovaquia",1360,35,355,988,8,2,17),
    @("Afganistan",1351,72,188,1120,7,1,43),
    @("Republica de Macedonia",1326,26,337,932,14,1,57),
    @("Cuba",1285,50,416,820,11,6,49),
    @("Ghana",1279,125,134,1135,4,1,10),
    @("Bulgaria",1188,91,193,941,37,2,54),
    @("Hong Kong",1036,0,699,333,9,0,4),
    @("Costa de Marfil",1004,0,359,631,0,0,14),
    @("Republica de Yibuti",999,13,330,667,0,0,2),
    @("Nigeria",981,0,197,753,2,0,31),
    @("Tunez",922,4,194,690,20,0,38),
    @("Guinea",862,0,170,686,0,0,6),
    @("Republica de Chipre",804,9,98,692,15,1,14),
    @("Letonia",784,6,267,505,6,1,12),
    @("Principado de Andorra",731,8,333,358,17,3,40),
    @("Crucero",712,0,645,54,4,0,13),
    @("Bolivia",703,31,44,616,3,3,43),
    @("Libano",696,8,140,534,46,0,22),
    @("Costa Rica",687,1,216,465,7,0,6),
    @("Albania",678,15,394,257,4,0,27),
    @("Niger",671,0,256,391,0,0,24),
    @("Kirguistan",656,25,322,326,11,0,8),
    @("Burkina Faso",616,0,410,165,0,0,41),
    @("Honduras",562,43,50,465,10,0,47),
    @("Uruguay",557,8,354,191,10,0,12),
    @("Senegal",545,66,262,276,1,1,7),
    @("San Marino",513,12,64,409,3,0,40),
    @("Estado de Palestina",484,4,92,388,0,0,4),
    @("Malta",447,2,223,221,2,0,3),
    @("Georgia",444,19,132,307,6,0,5),
    @("Jordania",441,4,326,108,5,0,7),
    @("Taiwan",428,1,264,158,0,0,6),
    @("Sri Lanka",417,49,109,301,2,0,7),
    @("Reunion",412,0,238,174,2,0,0),
    @("Consejo Danes para los Refugiados",394,17,48,321,0,0,25),
    @("Guatemala",384,42,30,343,5,1,11),
    @("Mayotte",354,28,144,206,4,0,4),
    @("Kenia",336,16,94,228,2,0,14),
    @("Mauricio",331,0,285,37,3,0,9),
    @("Somalia",328,0,8,304,2,0,16),
    @("Montenegro",319,3,123,190,7,1,6),
    @("Venezuela",318,20,132,176,4,0,10),
    @("Mali",309,0,77,211,0,0,21),
    @("Isla de Man",308,1,230,60,20,2,18),
    @("Tanzania",284,0,48,226,7,0,10),
    @("Vietnam",270,2,225,45,8,0,0),
    @("El Salvador",261,11,72,181,2,0,8),
    @("Jamaica",257,5,28,222,0,1,7),
    @("Paraguay",220,7,70,141,1,0,9),
    @("Guinea Ecuatorial",214,5,7,206,0,0,1),
    @("Islas Feroe",187,0,178,9,0,0,0),
    @("Congo",186,0,16,164,0,0,6),
    @("Sudan",174,12,14,144,0,3,16),
    @("Gabon",172,5,26,143,1,1,3),
    @("Martinica",170,6,77,79,6,0,14),
    @("Ruanda",154,0,87,67,0,0,0),
    @("Guadalupe",149,1,82,55,11,0,12),
    @("Birmania",144,12,9,130,0,0,5),
    @("Brunei",138,0,120,17,2,0,1),
    @("Gibraltar",133,0,129,4,0,0,0),
    @("Madagascar",122,1,61,61,1,0,0),
    @("Camboya",122,0,110,12,1,0,0),
    @("Etiopia",117,1,25,89,0,0,3),
    @("Liberia",117,16,25,84,0,0,8),
    @("Maldivas",116,8,16,100,2,0,0),
    @("Trinidad yTobago",115,0,48,59,0,0,8),
    @("Guayana Francesa",107,0,84,22,1,0,1),
    @("Aruba",100,0,68,30,4,0,2),
    @("Bermudas",99,0,39,55,10,0,5),
    @("Monaco",94,0,41,49,2,0,4),
    @("Togo",90,2,59,25,0,0,6),
    @("Cabo Verde",88,6,1,86,0,0,1),
    @("Zambia",84,8,37,44,1,0,3),
    @("Sierra Leona",82,18,10,70,0,1,2),
    @("Liechtenstein",81,0,55,25,0,0,1),
    @("Barbados",76,0,30,40,4,0,6),
    @("Uganda",75,1,46,29,0,0,0),
    @("Guyana",73,3,12,54,5,0,7),
    @("San Martin (Parte Holandesa)",73,0,22,39,8,0,12),
    @("Haiti",72,0,2,65,0,0,5),
    @("Bahamas",72,0,14,47,1,0,11),
    @("Islas Caimanes",66,0,7,58,3,0,1),
    @("Mozambique",65,19,12,53,0,0,0),
    @("Libia",60,0,18,40,0,0,2),
    @("Polinesia Francesa",57,0,36,21,1,0,0),
    @("Benin",54,0,27,26,0,0,1),
    @("Guinea-Bisau",50,0,3,47,0,0,0),
    @("Nepal",49,1,10,39,0,0,0),
    @("Macao",45,0,27,18,1,0,0),
    @("Siria",42,0,6,33,0,0,3),
    @("Republica del Chad",40,7,8,32,0,0,0),
    @("Puerto Rico",39,0,1,36,0,0,2),
    @("Eritrea",39,0,11,28,0,0,0),
    @("San Martin (Parte Francesa)",38,0,20,16,3,0,2),
    @("Mongolia",36,1,9,27,0,0,0),
    @("Suazilandia",36,5,10,25,0,0,1),
    @("Malaui",33,0,3,27,1,0,3),
    @("Guam",32,0,0,31,0,0,1),
    @("Zimbabue",29,1,2,23,0,0,4),
    @("Angola",25,0,6,17,0,0,2),
    @("Antigua y Barbuda",24,0,10,11,1,0,3),
    @("Timor Oriental",23,0,1,22,0,0,0),
    @("Botsuana",22,0,0,21,0,0,1),
    @("Laos",19,0,4,15,0,0,0),
    @("Belice",18,0,5,11,1,0,2),
    @("Fiyi",18,0,10,8,0,0,0),
    @("Nueva Caledonia",18,0,17,1,1,0,0),
    @("Islas Virgenes de los Estados Unidos",17,0,0,17,0,0,0),
    @("Namibia",16,0,7,9,0,0,0),
    @("Dominica",16,0,9,7,0,0,0),
    @("Republica de Africa Central",16,0,10,6,0,0,0),
    @("Curazao",16,2,11,4,0,0,1),
    @("San Cristobal y Nieves",15,0,2,13,0,0,0),
    @("Granada",15,0,7,8,4,0,0),
    @("Santa Lucia",15,0,15,0,0,0,0),
    @("San Vicente y las Granadinas",14,1,5,9,0,0,0),
    @("Islas Malvinas",12,0,11,1,0,0,0),
    @("Montserrat",11,0,2,9,1,0,0),
    @("Islas Turcas y Caicos",11,0,4,6,0,0,1),
    @("Burundi",11,0,4,6,0,0,1),
    @("Seychelles",11,0,6,5,0,0,0),
    @("Nicaragua",11,0,7,1,0,0,3),
    @("Groenlandia",11,0,11,0,0,0,0),
    @("Gambia",10,0,2,7,0,0,1),
    @("Surinam",10,0,6,3,0,0,1),
    @("Santa Sede",9,0,2,7,0,0,0),
    @("Papua Nueva Guinea",8,0,0,8,0,0,0),
    @("Butan",7,0,3,4,0,0,0),
    @("Mauritania",7,0,6,0,0,0,1),
    @("Sahara Occidental",6,0,5,1,0,0,0),
    @("San Bartolome",6,0,6,0,0,0,0),
    @("Bonaire, San Eustaquio y Saba",5,0,0,5,0,0,0),
    @("Sudan del Sur",5,1,0,5,0,0,0),
    @("Islas Virgenes Britanicas",5,0,3,1,0,0,1),
    @("Santo Tome y Principe",4,0,0,4,0,0,0),
    @("Anguila",3,0,1,2,0,0,0),
    @("San Pedro y Miquelon",1,0,0,1,0,0,0),
    @("Yemen",1,0,1,0,0,0,0)
)

$r = 4
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $r++
}
